$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 1245.4546
$ws.Range("J42").Value = 2163.3333
$ws.Range("L42").Value = 6489.999899999999
$ws.Range("N42").Value = -6949.999899999999
$ws.Range("H69").Value = 15112.25
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 18483
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 55449
$ws.Range("M69").Value = -14126
$ws.Range("N69").Value = -57197
$ws.Range("H72").Value = 15112.25
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 18483
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 166347
$ws.Range("M72").Value = -40632
$ws.Range("N72").Value = -175083
$ws.Range("H98").Value = 670.44446
$ws.Range("I98").Value = 679.3125
$ws.Range("J98").Value = 599.5
$ws.Range("K98").Value = 679.3125
$ws.Range("L98").Value = 599.5
$ws.Range("M98").Value = 818.6875
$ws.Range("N98").Value = -3595.5
$ws.Range("H100").Value = 1490
$ws.Range("I100").Value = 1042.1
$ws.Range("J100").Value = 1769.9375
$ws.Range("K100").Value = 1042.1
$ws.Range("L100").Value = 1769.9375
$ws.Range("M100").Value = -501.0999999999999
$ws.Range("N100").Value = -2851.9375
$ws.Range("H122").Value = 670.44446
$ws.Range("I122").Value = 679.3125
$ws.Range("J122").Value = 599.5
$ws.Range("K122").Value = 2037.9375
$ws.Range("L122").Value = 1798.5
$ws.Range("M122").Value = 412.0625
$ws.Range("N122").Value = -6698.5
$ws.Range("H132").Value = 9184.628000000001
$ws.Range("I132").Value = 6060.394
$ws.Range("K132").Value = 18181.182
$ws.Range("M132").Value = -15651.182

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4095.1667
$ws.Range("I61").Value = 3039.35
$ws.Range("K61").Value = 3039.35
$ws.Range("M61").Value = -2827.35
$ws.Range("H63").Value = 2730.4167
$ws.Range("J63").Value = 3538.3845
$ws.Range("L63").Value = 3538.3845
$ws.Range("N63").Value = -4910.3845
$ws.Range("H66").Value = 2730.4167
$ws.Range("J66").Value = 3538.3845
$ws.Range("L66").Value = 17691.9225
$ws.Range("N66").Value = -24555.9225
$ws.Range("H74").Value = 5640.413
$ws.Range("I74").Value = 5354.1816
$ws.Range("J74").Value = 6367
$ws.Range("K74").Value = 5354.1816
$ws.Range("L74").Value = 6367
$ws.Range("M74").Value = -4480.1816
$ws.Range("N74").Value = -8115
$ws.Range("H77").Value = 5640.413
$ws.Range("I77").Value = 5354.1816
$ws.Range("J77").Value = 6367
$ws.Range("K77").Value = 26770.908
$ws.Range("L77").Value = 31835
$ws.Range("M77").Value = -22402.908
$ws.Range("N77").Value = -40571
$ws.Range("H102").Value = 2972.889
$ws.Range("I102").Value = 2969.75
$ws.Range("K102").Value = 2969.75
$ws.Range("M102").Value = -1347.75
$ws.Range("H132").Value = 3008.5
$ws.Range("I132").Value = 2975.5151
$ws.Range("J132").Value = 3164
$ws.Range("K132").Value = 8926.5453
$ws.Range("L132").Value = 9492
$ws.Range("M132").Value = -6396.5453
$ws.Range("N132").Value = -14552
$ws.Range("H136").Value = 4095.1667
$ws.Range("I136").Value = 3039.35
$ws.Range("K136").Value = 9118.049999999999
$ws.Range("M136").Value = -6568.049999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2135.162
$ws.Range("I134").Value = 1972.25
$ws.Range("K134").Value = 5916.75
$ws.Range("M134").Value = -3381.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1554.762
$ws.Range("I22").Value = 1566.3636
$ws.Range("J22").Value = 1542
$ws.Range("K22").Value = 1566.3636
$ws.Range("L22").Value = 1542
$ws.Range("M22").Value = -1216.3636
$ws.Range("N22").Value = -2242
$ws.Range("H59").Value = 82000
$ws.Range("I59").Value = 48000
$ws.Range("J59").Value = 93333.336
$ws.Range("K59").Value = 48000
$ws.Range("L59").Value = 93333.336
$ws.Range("M59").Value = -46855
$ws.Range("N59").Value = -95623.336
$ws.Range("H122").Value = 3557.9783
$ws.Range("I122").Value = 2200.4443
$ws.Range("J122").Value = 4430.6787
$ws.Range("K122").Value = 6601.3329
$ws.Range("L122").Value = 13292.0361
$ws.Range("M122").Value = -4151.3329
$ws.Range("N122").Value = -18192.0361
$ws.Range("H134").Value = 6252.2036
$ws.Range("I134").Value = 6340.608
$ws.Range("K134").Value = 19021.824
$ws.Range("M134").Value = -16486.824

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 526450.1
$ws.Range("J9").Value = 202168
$ws.Range("L9").Value = 606504
$ws.Range("N9").Value = -606952
$ws.Range("H39").Value = 16282.167
$ws.Range("J39").Value = 19198.6
$ws.Range("L39").Value = 57595.8
$ws.Range("N39").Value = -58183.8
$ws.Range("H134").Value = 3798.5334
$ws.Range("I134").Value = 2994.5
$ws.Range("J134").Value = 5406.6
$ws.Range("K134").Value = 8983.5
$ws.Range("L134").Value = 16219.8
$ws.Range("M134").Value = -3913.5
$ws.Range("N134").Value = -26359.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5638.5557
$ws.Range("J132").Value = 5297
$ws.Range("L132").Value = 15891
$ws.Range("N132").Value = -20951

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I16").Value = 11365449
$ws.Range("J16").Value = 6999
$ws.Range("K16").Value = 11365449
$ws.Range("L16").Value = 6999
$ws.Range("M16").Value = -11365279
$ws.Range("N16").Value = -7339
$ws.Range("H132").Value = 63432.11
$ws.Range("I132").Value = 66810.47
$ws.Range("K132").Value = 200431.41
$ws.Range("M132").Value = -197901.41
$ws.Range("H136").Value = 5010679
$ws.Range("I136").Value = 12010308
$ws.Range("J136").Value = 10944.714
$ws.Range("K136").Value = 36030924
$ws.Range("L136").Value = 32834.142
$ws.Range("M136").Value = -36028374
$ws.Range("N136").Value = -37934.142

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 142864860
$ws.Range("I62").Value = 142864860
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 142864860
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -142864236
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 142864860
$ws.Range("I65").Value = 142864860
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 714324300
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -714321180
$ws.Range("N65").ClearContents()
$ws.Range("H81").Value = 5351521.5
$ws.Range("I81").Value = 5684741.5
$ws.Range("J81").Value = 20000
$ws.Range("K81").Value = 11369483
$ws.Range("L81").Value = 40000
$ws.Range("M81").Value = -11368422
$ws.Range("N81").Value = -42122
$ws.Range("H84").Value = 5351521.5
$ws.Range("I84").Value = 5684741.5
$ws.Range("J84").Value = 20000
$ws.Range("K84").Value = 56847415
$ws.Range("L84").Value = 200000
$ws.Range("M84").Value = -56842111
$ws.Range("N84").Value = -210608
$ws.Range("H100").Value = 853.6667
$ws.Range("I100").Value = 745.125
$ws.Range("J100").Value = 1070.75
$ws.Range("K100").Value = 1490.25
$ws.Range("L100").Value = 2141.5
$ws.Range("M100").Value = -949.25
$ws.Range("N100").Value = -3223.5
$ws.Range("H107").Value = 1026.7646
$ws.Range("I107").Value = 1566.5555
$ws.Range("J107").Value = 419.5
$ws.Range("K107").Value = 4699.666499999999
$ws.Range("L107").Value = 1258.5
$ws.Range("M107").Value = -2779.666499999999
$ws.Range("N107").Value = -5098.5
$ws.Range("H122").Value = 8613.226000000001
$ws.Range("I122").Value = 5765.8
$ws.Range("J122").Value = 20477.5
$ws.Range("K122").Value = 17297.4
$ws.Range("L122").Value = 61432.5
$ws.Range("M122").Value = -14847.4
$ws.Range("N122").Value = -66332.5
$ws.Range("H126").Value = 7332.4443
$ws.Range("I126").Value = 4040.0908
$ws.Range("J126").Value = 12506.143
$ws.Range("K126").Value = 12120.2724
$ws.Range("L126").Value = 37518.429
$ws.Range("M126").Value = -9650.2724
$ws.Range("N126").Value = -42458.429
$ws.Range("H132").Value = 4295.778
$ws.Range("I132").Value = 4545.375
$ws.Range("K132").Value = 13636.125
$ws.Range("M132").Value = -11106.125
$ws.Range("H136").Value = 2757
$ws.Range("I136").Value = 2649.4285
$ws.Range("K136").Value = 7948.2855
$ws.Range("M136").Value = -5398.2855
